# "Matching saledata w/ sn"
#
# The mall/erp_code/mall_name lookup table gets a 4th column ("mall_id",
# a numeric "sn" matched to each existing product row) and the erp_code
# header is renamed to erp_id.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- rename header --------------------------------------------------
$ws.Range("B1").Value = "erp_id"

# --- new column D: mall_id ------------------------------------------
# Row 2/6/7 end up styled "left + vertical-centered" and rows 1/3/4/5
# "left only" in the source file, so apply vertical-center *before* the
# left-only rows to reproduce that exact pair of cellXfs in that order.
$ws.Range("D2").Value = 45371
$ws.Range("D2").HorizontalAlignment = -4131   # xlLeft
$ws.Range("D2").VerticalAlignment = -4108     # xlCenter

$ws.Range("D1").Value = "mall_id"
$ws.Range("D1").HorizontalAlignment = -4131   # xlLeft

$ws.Range("D3").Value = 101716
$ws.Range("D3").HorizontalAlignment = -4131

$ws.Range("D4").Value = 104968
$ws.Range("D4").HorizontalAlignment = -4131

$ws.Range("D5").Value = 105250
$ws.Range("D5").HorizontalAlignment = -4131

$ws.Range("D6").Value = 104192
$ws.Range("D6").HorizontalAlignment = -4131
$ws.Range("D6").VerticalAlignment = -4108

$ws.Range("D7").Value = 108075
$ws.Range("D7").HorizontalAlignment = -4131
$ws.Range("D7").VerticalAlignment = -4108

# --- column widths ---------------------------------------------------
# Authored widths are 15.125 / 71.875 / 14.875 characters; this runtime's
# column grid only resolves to 1/7-character increments, so these inputs
# are chosen to land on the closest representable width.
$ws.Columns.Item(2).ColumnWidth = 14.428571428571429
$ws.Columns.Item(3).ColumnWidth = 71.14285714285714
$ws.Columns.Item(4).ColumnWidth = 14.142857142857142

# --- selection & print setup -----------------------------------------
$ws.Range("D5").Select()

$ws.PageSetup.PaperSize = 9      # xlPaperA4
$ws.PageSetup.Orientation = 1    # xlPortrait
